$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to text so numeric-looking strings
# (e.g. "587.82") are kept as exact text and not converted to floats.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "71.739.40"
$ws.Range("E2").Value = "  +3.37%  "
$ws.Range("D3").Value = "3.696.96"
$ws.Range("E3").Value = "  +7.11%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "587.82"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "180.70"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").Value = "3.685.64"
$ws.Range("E7").Value = "  +6.90%  "
$ws.Range("E8").Value = "  +3.74%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("E11").Value = "  +4.35%  "
$ws.Range("D12").Value = "49.85"
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "4.295.46"
$ws.Range("E14").Value = "  +7.33%  "
$ws.Range("D15").Value = "681.96"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").Value = "9.05"
$ws.Range("E16").Value = "  +3.99%  "
$ws.Range("D17").Value = "3.702.33"
$ws.Range("E17").Value = "  +7.32%  "
$ws.Range("D18").Value = "71.833.97"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").Value = "11.64"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").Value = "6.39"
$ws.Range("E22").Value = "  +18.92%  "
$ws.Range("E23").Value = "  +3.88%  "
$ws.Range("D24").Value = "17.81"
$ws.Range("E24").Value = "  +4.41%  "
$ws.Range("D25").Value = "103.41"
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("E26").Value = "  +3.26%  "
$ws.Range("D27").Value = "2.84"
$ws.Range("E27").Value = "  +4.79%  "
$ws.Range("E28").Value = "  +6.15%  "
$ws.Range("D29").Value = "35.60"
$ws.Range("E29").Value = "  +5.38%  "
$ws.Range("D30").Value = "9.32"
$ws.Range("E30").Value = "  +5.68%  "
$ws.Range("E31").Value = "  +6.54%  "
$ws.Range("D32").Value = "4.19"
$ws.Range("E32").Value = "  +11.78%  "
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").Value = "565.98"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").Value = "0.109"
$ws.Range("E35").Value = "  +3.62%  "
$ws.Range("D36").Value = "59.51"
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("D37").Value = "3.760.05"
$ws.Range("E37").Value = "  +3.44%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("D40").Value = "0.0₃0777"
$ws.Range("E40").Value = "  +4.53%  "
$ws.Range("D41").Value = "35.70"
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("E42").Value = "  +5.59%  "
$ws.Range("E43").Value = "  +4.19%  "
$ws.Range("E44").Value = "  +9.20%  "
$ws.Range("E45").Value = "  +4.65%  "
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").Value = "  +8.23%  "
$ws.Range("D47").Value = "3.38"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").Value = "0.133"
$ws.Range("E48").Value = "  +3.20%  "
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "135.86"
$ws.Range("E51").Value = "  +3.32%  "

# Restore the original (default) style for column D.
$ws.Range("D2:D51").Style = "Normal"
